$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): SMEs 1.6 -> 1.65, MSMEs 56.3 -> 56.35
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "1.65"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.35"

# Employment (% of total): SMEs 31.6 -> 31.63, MSMEs 73.2 -> 73.23
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "31.63"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "73.23"

# Enterprises (% of total): Micro 96.9 -> 96.95, SMEs 2.9 -> 2.92, MSMEs 99.9 -> 99.87
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "96.95"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "2.92"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "99.87"
